$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 251.5
$ws.Range("I39").Value = 235.33333
$ws.Range("K39").Value = 705.99999
$ws.Range("M39").Value = -409.99999

$ws.Range("H43").Value = 2437.8
$ws.Range("I43").Value = 2132
$ws.Range("J43").Value = 2896.5
$ws.Range("K43").Value = 2132
$ws.Range("L43").Value = 2896.5
$ws.Range("M43").Value = -2063
$ws.Range("N43").Value = -3034.5

$ws.Range("H80").Value = 198.63637
$ws.Range("I80").Value = 85
$ws.Range("J80").Value = 397.5
$ws.Range("K80").Value = 255
$ws.Range("L80").Value = 1192.5
$ws.Range("M80").Value = 743
$ws.Range("N80").Value = -3188.5

$ws.Range("H83").Value = 198.63637
$ws.Range("I83").Value = 85
$ws.Range("J83").Value = 397.5
$ws.Range("K83").Value = 765
$ws.Range("L83").Value = 3577.5
$ws.Range("M83").Value = 4227
$ws.Range("N83").Value = -13561.5

$ws.Range("H96").Value = 80.5
$ws.Range("I96").Value = 71
$ws.Range("K96").Value = 213
$ws.Range("M96").Value = 1160

$ws.Range("H100").Value = 738.6875
$ws.Range("I100").Value = 764.25
$ws.Range("K100").Value = 764.25
$ws.Range("M100").Value = -223.25

$ws.Range("H132").Value = 7486.394
$ws.Range("I132").Value = 7892.241
$ws.Range("J132").Value = 4544
$ws.Range("K132").Value = 23676.723
$ws.Range("L132").Value = 13632
$ws.Range("M132").Value = -21146.723
$ws.Range("N132").Value = -18692

$ws.Range("H141").Value = 992.7273
$ws.Range("I141").Value = 992.7273
$ws.Range("K141").Value = 2978.1819
$ws.Range("M141").Value = 2201.8181


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1081
$ws.Range("I2").Value = 1037
$ws.Range("J2").Value = 1345
$ws.Range("K2").Value = 1037
$ws.Range("L2").Value = 1345
$ws.Range("M2").Value = -924
$ws.Range("N2").Value = -1571

$ws.Range("H63").Value = 4177.222
$ws.Range("I63").Value = 2199.75
$ws.Range("K63").Value = 2199.75
$ws.Range("M63").Value = -1513.75

$ws.Range("H66").Value = 4177.222
$ws.Range("I66").Value = 2199.75
$ws.Range("K66").Value = 10998.75
$ws.Range("M66").Value = -7566.75

$ws.Range("H116").Value = 1081
$ws.Range("I116").Value = 1037
$ws.Range("J116").Value = 1345
$ws.Range("K116").Value = 1037
$ws.Range("L116").Value = 1345
$ws.Range("M116").Value = 1257
$ws.Range("N116").Value = -5933

$ws.Range("H122").Value = 2013.1666
$ws.Range("I122").Value = 2016.8
$ws.Range("K122").Value = 6050.4
$ws.Range("M122").Value = -3600.4


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1081
$ws.Range("I3").Value = 1037
$ws.Range("J3").Value = 1345
$ws.Range("K3").Value = 1037
$ws.Range("L3").Value = 1345
$ws.Range("M3").Value = -923
$ws.Range("N3").Value = -1573

$ws.Range("H94").Value = 724.5
$ws.Range("I94").Value = 724.5
$ws.Range("K94").Value = 724.5
$ws.Range("M94").Value = -273.5

$ws.Range("H99").Value = 1319.25
$ws.Range("I99").Value = 1319.25
$ws.Range("K99").Value = 1319.25
$ws.Range("M99").Value = 178.75

$ws.Range("H107").Value = 3478.8667
$ws.Range("I107").Value = 1168.3
$ws.Range("K107").Value = 1168.3
$ws.Range("M107").Value = 751.7


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6785.7896
$ws.Range("I31").Value = 2420.4
$ws.Range("J31").Value = 8344.857
$ws.Range("K31").Value = 2420.4
$ws.Range("L31").Value = 8344.857
$ws.Range("M31").Value = -2125.4
$ws.Range("N31").Value = -8934.857

$ws.Range("H34").Value = 6785.7896
$ws.Range("I34").Value = 2420.4
$ws.Range("J34").Value = 8344.857
$ws.Range("K34").Value = 2420.4
$ws.Range("L34").Value = 8344.857
$ws.Range("M34").Value = -2218.4
$ws.Range("N34").Value = -8748.857


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4158.3335
$ws.Range("I18").Value = 1633.3334
$ws.Range("K18").Value = 4900.0002
$ws.Range("M18").Value = -4731.0002

$ws.Range("H34").Value = 1914.0834
$ws.Range("I34").Value = 2638.25
$ws.Range("J34").Value = 1552
$ws.Range("K34").Value = 7914.75
$ws.Range("L34").Value = 4656
$ws.Range("M34").Value = -7830.75
$ws.Range("N34").Value = -4824

$ws.Range("H81").Value = 18997.5
$ws.Range("J81").Value = 18997.5
$ws.Range("L81").Value = 56992.5
$ws.Range("N81").Value = -59238.5

$ws.Range("H84").Value = 18997.5
$ws.Range("J84").Value = 18997.5
$ws.Range("L84").Value = 170977.5
$ws.Range("N84").Value = -182209.5

$ws.Range("H114").Value = 1474.1428
$ws.Range("I114").Value = 1981
$ws.Range("J114").Value = 1094
$ws.Range("K114").Value = 5943
$ws.Range("L114").Value = 3282
$ws.Range("M114").Value = -2689
$ws.Range("N114").Value = -9790


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 520.875
$ws.Range("I2").Value = 418
$ws.Range("J2").Value = 623.75
$ws.Range("K2").Value = 418
$ws.Range("L2").Value = 623.75
$ws.Range("M2").Value = -305
$ws.Range("N2").Value = -849.75

$ws.Range("H43").Value = 13407.647
$ws.Range("I43").Value = 1990
$ws.Range("J43").Value = 19635.455
$ws.Range("K43").Value = 1990
$ws.Range("L43").Value = 19635.455
$ws.Range("M43").Value = -1839
$ws.Range("N43").Value = -19937.455

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H122").Value = 2474.3333
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

$ws.Range("H132").Value = 2216
$ws.Range("I132").Value = 1873.6666
$ws.Range("K132").Value = 5620.9998
$ws.Range("M132").Value = -3090.9998

$ws.Range("H133").Value = 110000
$ws.Range("J133").Value = 110000
$ws.Range("L133").Value = 110000
$ws.Range("N133").Value = -120120

$ws.Range("H135").Value = 94000
$ws.Range("J135").Value = 94000
$ws.Range("L135").Value = 94000
$ws.Range("N135").Value = -104140


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6105.1113
$ws.Range("I40").Value = 4991.2856
$ws.Range("K40").Value = 4991.2856
$ws.Range("M40").Value = -4855.2856

$ws.Range("H122").Value = 2999.4
$ws.Range("I122").Value = 2874.25
$ws.Range("K122").Value = 8622.75
$ws.Range("M122").Value = -6172.75

$ws.Range("H136").Value = 1373.75
$ws.Range("I136").Value = 1373.75
$ws.Range("K136").Value = 4121.25
$ws.Range("M136").Value = -1571.25


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1173.5
$ws.Range("I113").Value = 1119.5714
$ws.Range("J113").Value = 1299.3334
$ws.Range("K113").Value = 3358.7142
$ws.Range("L113").Value = 3898.0002
$ws.Range("M113").Value = -1188.7142
$ws.Range("N113").Value = -8238.0002

$ws.Range("H122").Value = 1391.5
$ws.Range("I122").Value = 1216.5
$ws.Range("K122").Value = 3649.5
$ws.Range("M122").Value = -1199.5

$ws.Range("H132").Value = 1343.9667
$ws.Range("I132").Value = 1225.6786
$ws.Range("K132").Value = 3677.0358
$ws.Range("M132").Value = -1147.0358

